# Refresh cryptocurrency price / 1h-volume snapshot (GitHub Actions scrape).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '75.803.75'
$ws.Range("E2").Value = '  +9.40%  '
$ws.Range("D3").Value = '2.687.72'
$ws.Range("E3").Value = '  +11.13%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '''189.27'
$ws.Range("E5").Value = '  +14.26%  '
$ws.Range("D6").Value = '''590.33'
$ws.Range("E6").Value = '  +5.00%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("E8").Value = '  +5.76%  '
$ws.Range("D9").Value = '''0.198'
$ws.Range("E9").Value = '  +17.35%  '
$ws.Range("D10").Value = '2.688.72'
$ws.Range("E10").Value = '  +11.22%  '
$ws.Range("E11").Value = '  +1.59%  '
$ws.Range("D12").Value = '''0.360'
$ws.Range("E12").Value = '  +7.66%  '
$ws.Range("E13").Value = '  +1.83%  '
$ws.Range("D14").Value = '75.661.70'
$ws.Range("E14").Value = '  +9.46%  '
$ws.Range("D15").Value = '3.182.24'
$ws.Range("E15").Value = '  +11.12%  '
$ws.Range("E16").Value = '  +7.27%  '
$ws.Range("D17").Value = '''26.70'
$ws.Range("E17").Value = '  +11.57%  '
$ws.Range("D18").Value = '2.679.51'
$ws.Range("E18").Value = '  +10.73%  '
$ws.Range("D19").Value = '''9.44'
$ws.Range("E19").Value = '  +32.58%  '
$ws.Range("D20").Value = '''12.08'
$ws.Range("E20").Value = '  +11.88%  '
$ws.Range("D21").Value = '''375.07'
$ws.Range("E21").Value = '  +9.57%  '
$ws.Range("D22").Value = '''2.31'
$ws.Range("E22").Value = '  +18.26%  '
$ws.Range("D23").Value = '''4.07'
$ws.Range("E23").Value = '  +5.41%  '
$ws.Range("E24").Value = '  +4.46%  '
$ws.Range("D25").Value = '''1.00'
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("D26").Value = '''70.65'
$ws.Range("E26").Value = '  +7.09%  '
$ws.Range("D27").Value = '''4.20'
$ws.Range("E27").Value = '  +10.36%  '
$ws.Range("D28").Value = '''9.49'
$ws.Range("E28").Value = '  +11.91%  '
$ws.Range("D29").Value = '2.829.12'
$ws.Range("E29").Value = '  +11.15%  '
$ws.Range("E30").Value = '  +0.56%  '
$ws.Range("D31").Value = '0.0₃0963'
$ws.Range("E31").Value = '  +13.88%  '
$ws.Range("D32").Value = '''522.52'
$ws.Range("E32").Value = '  +15.19%  '
$ws.Range("E33").Value = '  +13.90%  '
$ws.Range("D34").Value = '''7.82'
$ws.Range("E34").Value = '  +6.03%  '
$ws.Range("E35").Value = '  +9.98%  '
$ws.Range("E36").Value = '  -0.07%  '
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").Value = '''0.120'
$ws.Range("E37").Value = '  +8.96%  '
$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").Value = '''162.50'
$ws.Range("E38").Value = '  +2.63%  '
$ws.Range("D39").Value = '''19.31'
$ws.Range("E39").Value = '  +6.29%  '
$ws.Range("D40").Value = '''19.38'
$ws.Range("E40").Value = '  +1.53%  '
$ws.Range("E41").Value = '  +0.07%  '
$ws.Range("E42").Value = '  +15.14%  '
$ws.Range("D43").Value = '''171.16'
$ws.Range("E43").Value = '  +26.79%  '
$ws.Range("D44").Value = '''1.71'
$ws.Range("E44").Value = '  +12.69%  '
$ws.Range("E45").Value = '  +10.49%  '
$ws.Range("E46").Value = '  +11.16%  '
$ws.Range("E47").Value = '  +16.03%  '
$ws.Range("D48").Value = '''39.30'
$ws.Range("E48").Value = '  +4.06%  '
$ws.Range("D49").Value = '''0.0848'
$ws.Range("E49").Value = '  +17.26%  '
$ws.Range("E50").Value = '  +8.49%  '
$ws.Range("E51").Value = '  +10.97%  '
